$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores Price/Volume columns as literal text (inline strings),
# preserving formatting such as thousands separators written as extra dots and
# padded percent signs. Where the new Price text would otherwise be auto-
# detected by Excel as a plain number, force the cell to Text format first so
# the literal characters are preserved on save.

$ws.Range("D2").Value = '71.967.84'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '4.009.28'
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.83'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.17'
$ws.Range("E6").Value = '  +2.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.696'
$ws.Range("E7").Value = '  +10.61%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  -3.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000326'
$ws.Range("E11").Value = '  -5.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.66'
$ws.Range("E12").Value = '  +5.96%  '
$ws.Range("E13").Value = '  -1.47%  '
$ws.Range("D14").Value = '4.651.43'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '4.014.94'
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.95'
$ws.Range("E16").Value = '  -3.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.51'
$ws.Range("E17").Value = '  -4.26%  '
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("D20").Value = '71.828.03'
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '425.99'
$ws.Range("E21").Value = '  -4.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '97.78'
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.49'
$ws.Range("E23").Value = '  -3.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.19'
$ws.Range("E24").Value = '  +1.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.28'
$ws.Range("E26").Value = '  -9.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.72'
$ws.Range("E27").Value = '  -3.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.84'
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.68'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.56'
$ws.Range("E30").Value = '  +22.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.34'
$ws.Range("E31").Value = '  -3.16%  '
$ws.Range("E32").Value = '  -1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.12'
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '675.14'
$ws.Range("E34").Value = '  -4.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.54'
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '43.96'
$ws.Range("E36").Value = '  +5.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.434'
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.153'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").Value = '0.0₃0827'
$ws.Range("E39").Value = '  -9.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.44'
$ws.Range("E40").Value = '  -4.14%  '
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.26'
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("E48").Value = '  -8.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.98'
$ws.Range("E49").Value = '  -7.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000273'
$ws.Range("E50").Value = '  -3.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '145.00'
$ws.Range("E51").Value = '  +1.18%  '

# Rows 46 and 47 swap content (THORChain <-> ApeXProtocol traded ranks)
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.42'
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.63'
$ws.Range("E47").Value = '  +2.74%  '
